$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of C20:F20 (keep existing cell formatting/style)
$ws.Range("C20:F20").ClearContents()

# Scroll the view so row 10 is at the top, then select cell I20,
# matching the saved view state (topLeftCell="A10", selection I20)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I20").Select()
